$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header row (A1:D1) to snake_case column names
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# 2. Title-case connector words (de/del/la/las/los/el/y) in state/municipality names
#    plus two special case fixes (GUANAJUATO -> Guanajuato, MonteMorelos -> Montemorelos)
$ws.Range('B8').Value = 'Pabellón De Arteaga'
$ws.Range('B9').Value = 'Rincón De Romos'
$ws.Range('B10').Value = 'San Francisco De Los Romo'
$ws.Range('B11').Value = 'San José De Gracia'
$ws.Range('B16').Value = 'Playas De Rosarito'
$ws.Range('B40').Value = 'Amatenango De La Frontera'
$ws.Range('B41').Value = 'Amatenango Del Valle'
$ws.Range('B44').Value = 'Bejucal De Ocampo'
$ws.Range('B46').Value = 'Benemérito De Las Américas'
$ws.Range('B53').Value = 'Chiapa De Corzo'
$ws.Range('B58').Value = 'Comitán De Domínguez'
$ws.Range('B78').Value = 'Mazapa De Madero'
$ws.Range('B83').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B91').Value = 'Salto De Agua'
$ws.Range('B92').Value = 'San Cristóbal De Las Casas'
$ws.Range('B129').Value = 'Coyame Del Sotol'
$ws.Range('B140').Value = 'Guadalupe Y Calvo'
$ws.Range('B143').Value = 'Hidalgo Del Parral'
$ws.Range('B167').Value = 'San Francisco De Borja'
$ws.Range('B168').Value = 'San Francisco De Conchos'
$ws.Range('B169').Value = 'San Francisco Del Oro'
$ws.Range('B177').Value = 'Valle De Zaragoza'
$ws.Range('B210').Value = 'San Juan De Sabinas'
$ws.Range('B227').Value = 'Villa De Álvarez'
$ws.Range('A229').Value = 'Ciudad De México'
$ws.Range('B233').Value = 'Cuajimalpa De Morelos'
$ws.Range('B248').Value = 'Coneto De Comonfort'
$ws.Range('B262').Value = 'Nombre De Dios'
$ws.Range('B266').Value = 'Pánuco De Coronado'
$ws.Range('B273').Value = 'San Juan De Guadalupe'
$ws.Range('B274').Value = 'San Juan Del Río'
$ws.Range('B275').Value = 'San Luis Del Cordero'
$ws.Range('B276').Value = 'San Pedro Del Gallo'
$ws.Range('A286').Value = 'Estado De México'
$ws.Range('B286').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B289').Value = 'Almoloya De Alquisiras'
$ws.Range('B290').Value = 'Almoloya De Juárez'
$ws.Range('B297').Value = 'Atizapán De Zaragoza'
$ws.Range('B304').Value = 'Chapa De Mota'
$ws.Range('B309').Value = 'Coacalco De Berriozábal'
$ws.Range('B316').Value = 'Ecatepec De Morelos'
$ws.Range('B324').Value = 'Ixtapan De La Sal'
$ws.Range('B325').Value = 'Ixtapan Del Oro'
$ws.Range('B340').Value = 'Naucalpan De Juárez'
$ws.Range('B352').Value = 'San Antonio La Isla'
$ws.Range('B353').Value = 'San Felipe Del Progreso'
$ws.Range('B354').Value = 'San Martín De Las Pirámides'
$ws.Range('B356').Value = 'San Simón De Guerero'
$ws.Range('B358').Value = 'Soyaniquilpan De Juárez'
$ws.Range('B367').Value = 'Tenango Del Valle'
$ws.Range('B381').Value = 'Tlalnepantla De Baz'
$ws.Range('B387').Value = 'Valle De Bravo'
$ws.Range('B388').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B389').Value = 'Villa De Allende'
$ws.Range('B390').Value = 'Villa Del Carbón'
$ws.Range('B404').Value = 'Apaseo El Alto'
$ws.Range('B405').Value = 'Apaseo El Grande'
$ws.Range('B413').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B417').Value = 'Jaral Del Progreso'
$ws.Range('B425').Value = 'Purísima Del Rincón'
$ws.Range('B429').Value = 'San Diego De La Unión'
$ws.Range('B431').Value = 'San Francisco Del Rincón'
$ws.Range('B433').Value = 'San Luis De La Paz'
$ws.Range('B435').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B437').Value = 'Silao De La Victoria'
$ws.Range('B442').Value = 'Valle De Santiago'
$ws.Range('B448').Value = 'Acapulco De Juárez'
$ws.Range('B451').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B452').Value = 'Alcozauca De Guerero'
$ws.Range('B456').Value = 'Atenango Del Río'
$ws.Range('B458').Value = 'Atoyac De Álvarez'
$ws.Range('B459').Value = 'Ayutla De Los Libres'
$ws.Range('B462').Value = 'Buenavista De Cuéllar'
$ws.Range('B463').Value = 'Chilapa De Álvarez'
$ws.Range('B464').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B465').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B470').Value = 'Coyuca De Benítez'
$ws.Range('B471').Value = 'Coyuca De Catalán'
$ws.Range('B475').Value = 'Cuetzala Del Progreso'
$ws.Range('B476').Value = 'Cutzamala De Pinzón'
$ws.Range('B482').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B483').Value = 'Iguala De La Independencia'
$ws.Range('B485').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B486').Value = 'Zihuatanejo De Azueta'
$ws.Range('B488').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B491').Value = 'Mártir De Cuilapan'
$ws.Range('B504').Value = 'Taxco De Alarcón'
$ws.Range('B506').Value = 'Técpan De Galeana'
$ws.Range('B508').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B510').Value = 'Tixtla De Guerero'
$ws.Range('B514').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B515').Value = 'Tlapa De Comonfort'
$ws.Range('B527').Value = 'Agua Blanca De Iturbide'
$ws.Range('B534').Value = 'Atotonilco De Tula'
$ws.Range('B535').Value = 'Atotonilco El Grande'
$ws.Range('B541').Value = 'Cuautepec De Hinojosa'
$ws.Range('B547').Value = 'Huasca De Ocampo'
$ws.Range('B550').Value = 'Huejutla De Reyes'
$ws.Range('B553').Value = 'Jacala De Ledezma'
$ws.Range('B559').Value = 'Mineral Del Chico'
$ws.Range('B560').Value = 'Mineral Del Monte'
$ws.Range('B561').Value = 'Mixquiahuala De Juárez'
$ws.Range('B562').Value = 'Molango De Escamilla'
$ws.Range('B564').Value = 'Nopala De Villagrán'
$ws.Range('B565').Value = 'Omitlán De Juárez'
$ws.Range('B566').Value = 'Pachuca De Soto'
$ws.Range('B569').Value = 'Progreso De Obregón'
$ws.Range('B574').Value = 'Santiago De Anaya'
$ws.Range('B575').Value = 'Santiago Tulantepec De Lugo Guerero'
$ws.Range('B579').Value = 'Tenango De Doria'
$ws.Range('B581').Value = 'Tepehuacán De Guerero'
$ws.Range('B582').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B584').Value = 'Tezontepec De Aldama'
$ws.Range('B592').Value = 'Tula De Allende'
$ws.Range('B593').Value = 'Tulancingo De Bravo'
$ws.Range('B594').Value = 'Villa De Tezontepec'
$ws.Range('B597').Value = 'Zacualtipán De Ángeles'
$ws.Range('B598').Value = 'Zapotlán De Juárez'
$ws.Range('B603').Value = 'Acatlán De Juárez'
$ws.Range('B604').Value = 'Ahualulco De Mercado'
$ws.Range('B608').Value = 'Atemajac De Brizuela'
$ws.Range('B610').Value = 'Atotonilco El Alto'
$ws.Range('B612').Value = 'Autlán De Navarro'
$ws.Range('B618').Value = 'Cañadas De Obregón'
$ws.Range('B625').Value = 'Concepción De Buenos Aires'
$ws.Range('B626').Value = 'Cuautitlán De García Barragán'
$ws.Range('B633').Value = 'Encarnación De Díaz'
$ws.Range('B639').Value = 'Huejuquilla El Alto'
$ws.Range('B640').Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range('B641').Value = 'Ixtlahuacán Del Río'
$ws.Range('B645').Value = 'Jilotlán De Los Dolores'
$ws.Range('B651').Value = 'La Manzanilla De La Paz'
$ws.Range('B652').Value = 'Lagos De Moreno'
$ws.Range('B660').Value = 'Ojuelos De Jalisco'
$ws.Range('B665').Value = 'San Cristóbal De La Barranca'
$ws.Range('B666').Value = 'San Diego De Alejandría'
$ws.Range('B667').Value = 'San Juan De Los Lagos'
$ws.Range('B669').Value = 'San Martín De Bolaños'
$ws.Range('B671').Value = 'San Miguel El Alto'
$ws.Range('B672').Value = 'San Sebastián Del Oeste'
$ws.Range('B673').Value = 'Santa María De Los Ángeles'
$ws.Range('B674').Value = 'Santa María Del Oro'
$ws.Range('B677').Value = 'Talpa De Allende'
$ws.Range('B678').Value = 'Tamazula De Gordiano'
$ws.Range('B681').Value = 'Techaluta De Montenegro'
$ws.Range('B685').Value = 'Teocuitatlán De Corona'
$ws.Range('B686').Value = 'Tepatitlán De Morelos'
$ws.Range('B689').Value = 'Tizapán El Alto'
$ws.Range('B690').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B701').Value = 'Unión De San Antonio'
$ws.Range('B702').Value = 'Unión De Tula'
$ws.Range('B703').Value = 'Valle De Guadalupe'
$ws.Range('B704').Value = 'Valle De Juárez'
$ws.Range('B709').Value = 'Yahualica De González Gallo'
$ws.Range('B710').Value = 'Zacoalco De Torres'
$ws.Range('B713').Value = 'Zapotlán Del Rey'
$ws.Range('B714').Value = 'Zapotlán El Grande'
$ws.Range('B738').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B740').Value = 'Cojumatlán De Régules'
$ws.Range('B806').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B832').Value = 'Coatlán Del Río'
$ws.Range('B839').Value = 'Jonacatepec De Leandro Valle'
$ws.Range('B843').Value = 'Puente De Ixtla'
$ws.Range('B849').Value = 'Tetela Del Volcán'
$ws.Range('B850').Value = 'Tlaltizapán De Zapata'
$ws.Range('B857').Value = 'Zacualpan De Amilpas'
$ws.Range('B861').Value = 'Amatlán De Cañas'
$ws.Range('B862').Value = 'Bahía De Banderas'
$ws.Range('B866').Value = 'Ixtlán Del Río'
$ws.Range('B872').Value = 'Santa María Del Oro'
$ws.Range('B890').Value = 'Ciénega De Flores'
$ws.Range('B908').Value = 'Lampazos De Naranjo'
$ws.Range('B915').Value = 'Mier Y Noriega'
$ws.Range('B924').Value = 'San Nicolás De Los Garza'
$ws.Range('B931').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B936').Value = 'Ciénega De Zimatlán'
$ws.Range('B938').Value = 'Coicoyán De Las Flores'
$ws.Range('B941').Value = 'Cuilápam De Guerero'
$ws.Range('B942').Value = 'Fresnillo De Trujano'
$ws.Range('B943').Value = 'Guadalupe De Ramírez'
$ws.Range('B944').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B945').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B946').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B947').Value = 'Huautla De Jiménez'
$ws.Range('B948').Value = 'Ixtlán De Juárez'
$ws.Range('B949').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B955').Value = 'Mariscala De Juárez'
$ws.Range('B957').Value = 'Mazatlán Villa De Flores'
$ws.Range('B958').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B960').Value = 'Nejapa De Madero'
$ws.Range('B961').Value = 'Oaxaca De Juárez'
$ws.Range('B962').Value = 'Ocotlán De Morelos'
$ws.Range('B963').Value = 'Putla Villa De Guerero'
$ws.Range('B964').Value = 'Reforma De Pineda'
$ws.Range('B973').Value = 'San Antonino El Alto'
$ws.Range('B982').Value = 'San Felipe Jalapa De Díaz'
$ws.Range('B986').Value = 'San Francisco Del Mar'
$ws.Range('B997').Value = 'San José Del Progreso'
$ws.Range('B1004').Value = 'San Juan De Los Cués'
$ws.Range('B1005').Value = 'San Juan Del Río'
$ws.Range('B1015').Value = 'San Martín De Los Cansecos'
$ws.Range('B1019').Value = 'San Miguel Del Puerto'
$ws.Range('B1050').Value = 'Santa Ana Del Valle'
$ws.Range('B1059').Value = 'Santa Cruz De Bravo'
$ws.Range('B1067').Value = 'Santa Inés Del Monte'
$ws.Range('B1077').Value = 'Santa María Jalapa Del Marqués'
$ws.Range('B1117').Value = 'Santo Domingo De Morelos'
$ws.Range('B1134').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B1136').Value = 'Tataltepec De Valdés'
$ws.Range('B1137').Value = 'Teotitlán De Flores Magón'
$ws.Range('B1139').Value = 'Tepelmeme Villa De Morelos'
$ws.Range('B1140').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B1141').Value = 'Tlacolula De Matamoros'
$ws.Range('B1143').Value = 'Totontepec Villa De Morelos'
$ws.Range('B1146').Value = 'Villa De Chilapa De Díaz'
$ws.Range('B1147').Value = 'Villa De Etla'
$ws.Range('B1148').Value = 'Villa De Tamazulápam Del Progreso'
$ws.Range('B1149').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B1150').Value = 'Villa De Zaachila'
$ws.Range('B1152').Value = 'Villa Sola De Vega'
$ws.Range('B1153').Value = 'Villa Talea De Castro'
$ws.Range('B1154').Value = 'Zapotitlán Del Río'
$ws.Range('B1156').Value = 'Zimatlán De Álvarez'
$ws.Range('B1176').Value = 'Ayotoxco De Guerero'
$ws.Range('B1180').Value = 'Chalchicomula De Sesma'
$ws.Range('B1191').Value = 'Chila De La Sal'
$ws.Range('B1198').Value = 'Cuayuca De Andrade'
$ws.Range('B1199').Value = 'Cuetzalan Del Progreso'
$ws.Range('B1213').Value = 'Huehuetlán El Chico'
$ws.Range('B1214').Value = 'Huehuetlán El Grande'
$ws.Range('B1218').Value = 'Huitzilan De Serdán'
$ws.Range('B1220').Value = 'Ixcamilpa De Guerero'
$ws.Range('B1224').Value = 'Izúcar De Matamoros'
$ws.Range('B1232').Value = 'Los Reyes De Juárez'
$ws.Range('B1242').Value = 'Palmar De Bravo'
$ws.Range('B1260').Value = 'San Salvador El Seco'
$ws.Range('B1261').Value = 'San Salvador El Verde'
$ws.Range('B1268').Value = 'Tecali De Herrera'
$ws.Range('B1276').Value = 'Tepanco De López'
$ws.Range('B1277').Value = 'Tepango De Rodríguez'
$ws.Range('B1278').Value = 'Tepatlaxco De Hidalgo'
$ws.Range('B1284').Value = 'Tepexi De Rodríguez'
$ws.Range('B1286').Value = 'Tepeyahualco De Cuauhtémoc'
$ws.Range('B1287').Value = 'Tetela De Ocampo'
$ws.Range('B1288').Value = 'Teteles De Avila Castillo'
$ws.Range('B1293').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B1305').Value = 'Totoltepec De Guerero'
$ws.Range('B1307').Value = 'Tuzamapan De Galeana'
$ws.Range('B1322').Value = 'Zapotitlán De Méndez'
$ws.Range('B1327').Value = 'Amealco De Bonfil'
$ws.Range('B1329').Value = 'Cadereyta De Montes'
$ws.Range('B1335').Value = 'Jalpan De Serra'
$ws.Range('B1336').Value = 'Landa De Matamoros'
$ws.Range('B1339').Value = 'Pinal De Amoles'
$ws.Range('B1342').Value = 'San Juan Del Río'
$ws.Range('B1357').Value = 'Armadillo De Los Infante'
$ws.Range('B1358').Value = 'Axtla De Terrazas'
$ws.Range('B1363').Value = 'Cerro De San Pedro'
$ws.Range('B1365').Value = 'Ciudad Del Maíz'
$ws.Range('B1376').Value = 'Mexquitic De Carmona'
$ws.Range('B1382').Value = 'San Ciro De Acosta'
$ws.Range('B1388').Value = 'Santa María Del Río'
$ws.Range('B1390').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B1398').Value = 'Tanquián De Escobedo'
$ws.Range('B1402').Value = 'Villa De Arista'
$ws.Range('B1403').Value = 'Villa De Arriaga'
$ws.Range('B1404').Value = 'Villa De Guadalupe'
$ws.Range('B1405').Value = 'Villa De La Paz'
$ws.Range('B1406').Value = 'Villa De Ramos'
$ws.Range('B1407').Value = 'Villa De Reyes'
$ws.Range('B1448').Value = 'Nacozari De García'
$ws.Range('B1466').Value = 'Jalpa De Méndez'
$ws.Range('B1511').Value = 'Soto La Marina'
$ws.Range('B1519').Value = 'Acuamanala De Miguel Hidalgo'
$ws.Range('B1521').Value = 'Amaxac De Guerero'
$ws.Range('B1525').Value = 'Contla De Juan Cuamatzi'
$ws.Range('B1530').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B1533').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B1536').Value = 'Papalotla De Xicohténcatl'
$ws.Range('B1538').Value = 'San Pablo Del Monte'
$ws.Range('B1542').Value = 'Tepetitla De Lardizábal'
$ws.Range('B1545').Value = 'Tetla De La Solidaridad'
$ws.Range('B1562').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B1566').Value = 'Amatlán De Los Reyes'
$ws.Range('B1573').Value = 'Boca Del Río'
$ws.Range('B1575').Value = 'Camarón De Tejeda'
$ws.Range('B1579').Value = 'Castillo De Teayo'
$ws.Range('B1581').Value = 'Cazones De Herrera'
$ws.Range('B1589').Value = 'Chinampa De Gorostiza'
$ws.Range('B1602').Value = 'Cosamaloapan De Carpio'
$ws.Range('B1603').Value = 'Cosautlán De Carvajal'
$ws.Range('B1620').Value = 'Hueyapan De Ocampo'
$ws.Range('B1621').Value = 'Ignacio De La Llave'
$ws.Range('B1625').Value = 'Ixhuatlán De Madero'
$ws.Range('B1626').Value = 'Ixhuatlán Del Café'
$ws.Range('B1627').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B1638').Value = 'Juchique De Ferrer'
$ws.Range('B1640').Value = 'Landero Y Coss'
$ws.Range('B1642').Value = 'Las Vigas De Ramírez'
$ws.Range('B1643').Value = 'Lerdo De Tejada'
$ws.Range('B1648').Value = 'Martínez De La Torre'
$ws.Range('B1651').Value = 'Medellín De Bravo'
$ws.Range('B1655').Value = 'Mixtla De Altamirano'
$ws.Range('B1657').Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range('B1667').Value = 'Ozuluama De Mascareñas'
$ws.Range('B1671').Value = 'Paso De Ovejas'
$ws.Range('B1672').Value = 'Paso Del Macho'
$ws.Range('B1676').Value = 'Poza Rica De Hidalgo'
$ws.Range('B1686').Value = 'Sayula De Alemán'
$ws.Range('B1690').Value = 'Soledad De Doblado'
$ws.Range('B1698').Value = 'Tatahuicapan De Juárez'
$ws.Range('B1721').Value = 'Tlacotepec De Mejía'
$ws.Range('B1735').Value = 'Vega De Alatorre'
$ws.Range('B1745').Value = 'Zontecomatlán De López Y Fuentes'
$ws.Range('B1746').Value = 'Zozocolco De Hidalgo'
$ws.Range('B1780').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B1782').Value = 'Concepción Del Oro'
$ws.Range('B1794').Value = 'Jiménez Del Teul'
$ws.Range('B1801').Value = 'Mezquital Del Oro'
$ws.Range('B1806').Value = 'Moyahua De Estrada'
$ws.Range('B1807').Value = 'Nochistlán De Mejía'
$ws.Range('B1808').Value = 'Noria De Ángeles'
$ws.Range('B1819').Value = 'Teúl De González Ortega'
$ws.Range('B1820').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B1822').Value = 'Trinidad García De La Cadena'
$ws.Range('B1825').Value = 'Villa De Cos'
$ws.Range('A401').Value = 'Guanajuato'
$ws.Range('B917').Value = 'Montemorelos'

# 3. Fix floating point precision drift in percentage column (recalculation artifacts)
$ws.Range('D268').Value = 0.0009201592359772576
$ws.Range('D401').Value = 0.0009589027827552472
$ws.Range('D821').Value = 0.00094921689606075
$ws.Range('D1336').Value = 0.00091047334928276
$ws.Range('D1602').Value = 0.00091047334928276
$ws.Range('D1738').Value = 0.0009589027827552472

# 4. Remove trailing metadata/footer rows 1834-1838 (sample size, source, author, date notes)
$ws.Range("A1834:D1838").EntireRow.Delete()

# Dimension will auto-update to A1:D1832 after the row deletion